$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DDL strings for the create-table text that goes into column B.
$sqlAppInstanceData = @"
create table APP_INSTANCE_DATA
(
     APP_ID SERIAL PRIMARY KEY,
     APP_NAME CHARACTER VARYING(255) NOT NULL,
     APP_TYPE CHARACTER VARYING(255) NOT NULL
);
"@

$sqlAppInstancePayload = @"
create table APP_INSTANCE_PAYLOAD
(   
    APP_ID INTEGER NOT NULL,
    PAYLOAD BYTEA NOT NULL
);

"@

$sqlTeamApps = @"
create table TEAM_APPS
(
APP_ID INTEGER NOT NULL,
TEAM_ID INTEGER NOT NULL
);

"@

$sqlTeams = @"
create table TEAMS
(
    TEAM_ID INTEGER NOT NULL,
    TEAM_NAME CHARACTER VARYING(255) NOT NULL,
    TEAM_DL CHARACTER VARYING(255) NOT NULL
);
"@

$sqlUserTeamRelation = @"

create table USER_TEAM_RELATION
(
USER_ID INTEGER NOT NULL,
TEAM_ID INTEGER NOT NULL
);

"@

$sqlPersonalApps = @"
create table PERSONAL_APPS
(
USER_ID INTEGER NOT NULL,
APP_ID INTEGER NOT NULL
);
"@

$sqlGlobalApps = @"
create table GLOBAL_APPS
(
APP_ID INTEGER NOT NULL
);
"@

# Add a second column to hold the DDL text and set its width (~61.13 chars).
$ws.Columns.Item(2).ColumnWidth = 60.3

# --- Section 1: APP_INSTANCE_DATA (rows 1-4) ---
$ws.Range("B1").Value = $sqlAppInstanceData
$ws.Range("B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 79.85

# --- Section 2: APP_INSTANCE_PAYLOAD (rows 7-9) ---
$ws.Range("B7").Value = $sqlAppInstancePayload
$ws.Range("B7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 68.65

# --- Section 3: TEAM_APPS (rows 11-13) ---
$ws.Range("B11").Value = $sqlTeamApps
$ws.Range("B11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 68.65

# --- Section 4: TEAMS (rows 15-18) ---
$ws.Range("B15").Value = $sqlTeams
$ws.Range("B15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 68.65

# --- Section 5: USER_TEAM_RELATION (rows 20-22) ---
$ws.Range("B20").Value = $sqlUserTeamRelation
$ws.Range("B20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 79.85

# --- Section 6: PERSONAL_APPS (rows 24-26) ---
$ws.Range("B24").Value = $sqlPersonalApps
$ws.Range("B24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 57.45

# --- Section 7: GLOBAL_APPS (rows 28) ---
$ws.Range("B28").Value = $sqlGlobalApps
$ws.Range("B28").WrapText = $true
$ws.Rows.Item(28).RowHeight = 46.25

# Update selection / view to match the target state
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("B29").Select()
